# Applies the "Updates to U.S. develop commit #d300b00a" change:
#  - Adds two new note rows (15 and 16) to the "About" sheet describing the
#    assumption that new coal with CCS cannot be built before 2028.
#  - On the "BBNPPTY" sheet, flips the Boolean flags for years 2024-2027
#    (columns E:H) from 0 to 1 for the four "...w CCS" technology rows
#    (hard coal w CCS, natural gas combined cycle w CCS, biomass w CCS,
#    lignite w CCS -- rows 19-22).
#  - Leaves "About" as the active/selected sheet with the last selection on
#    B22, and leaves the BBNPPTY selection on E19:H22.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsBBNPPTY = $wb.Worksheets.Item("BBNPPTY")

# --- Add the two new explanatory note rows under the existing note on the
# --- About sheet (rows 14 is blank, so the new text lands on 15 and 16).
$wsAbout.Range("A15").Value = "We also assume no new coal with CCS can be built prior to 2028 given the state of the technology"
$wsAbout.Range("A16").Value = "and the construction time for new or modified plants."

# --- Flip the 2024-2027 (columns E:H) Boolean ban flags to 1 for the four
# --- "w CCS" rows on the BBNPPTY sheet.
$wsBBNPPTY.Range("E19:H22").Value = 1

# --- Restore the selections/active sheet to match the saved view state.
$null = $wsBBNPPTY.Range("E19:H22").Select()
$null = $wsAbout.Activate()
$null = $wsAbout.Range("B22").Select()
